# Config.xlsx update: robot now pulls transaction details from the
# Orchestrator Queue instead of local/IFL specific constants, and the
# Assets sheet gains a FolderLocation_ExecutionFile asset entry.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Settings"
$ws3 = $wb.Worksheets.Item(3)   # "Assets"

# --- Settings sheet -------------------------------------------------
# Row 2 (OrchestratorQueueName): value renamed from TestExecutionQueue
# to ExecutionQueue.
# Row 3 (OrchestratorFolderPath): value changed from the old IFL folder
# name to "Generic Asset"; its description cell is cleared.
# NOTE: set B3 before B2 so the shared-string table is built in the
# same order the author's file uses ("Generic Asset" before
# "ExecutionQueue").
$ws1.Range("B3").Value = "Generic Asset"
$ws1.Range("B2").Value = "ExecutionQueue"
$ws1.Range("C3").ClearContents()

# Rows 4-6 (SharedRepository, DelayShort, DelayLong) are removed
# entirely - they are no longer part of the config. This shifts the
# blank spacer row and the logF_BusinessProcessName row up by three.
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(4).Delete()
$ws1.Rows.Item(4).Delete()

# restore selection to match the saved workbook state
[void]$ws1.Activate()
[void]$ws1.Range("A3").Select()

# --- Assets sheet -----------------------------------------------------
# New asset row: FolderLocation_ExecutionFile, pointing at the
# "Generic Asset" orchestrator asset folder.
$ws3.Range("A2").Value = "FolderLocation_ExecutionFile"
$ws3.Range("B2").Value = "FolderLocation_ExecutionFile"
$ws3.Range("C2").Value = "Generic Asset"

# the sheet's trailing blank row (the used range shrank by one row
# once the sheet was touched/resaved)
$ws3.Rows.Item(1000).Delete()

[void]$ws3.Activate()
[void]$ws3.Range("A8").Select()

# leave "Settings" as the active tab, matching the original workbook
[void]$ws1.Activate()
